# Update "想去人数" (want-to-attend count) values in column F across sheets,
# matching the scraper re-run captured in the diff.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 37563
$ws.Range("F5").Value = 767
$ws.Range("F6").Value = 475
$ws.Range("F9").Value = 841
$ws.Range("F11").Value = 708
$ws.Range("F13").Value = 34
$ws.Range("F15").Value = 14
$ws.Range("F16").Value = 646
$ws.Range("F20").Value = 1164
$ws.Range("F22").Value = 822
$ws.Range("F23").Value = 2514
$ws.Range("F24").Value = 997
$ws.Range("F25").Value = 561
$ws.Range("F26").Value = 106
$ws.Range("F27").Value = 1156
$ws.Range("F29").Value = 765
$ws.Range("F30").Value = 55
$ws.Range("F31").Value = 1152

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 384
$ws.Range("F11").Value = 9

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 628

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 628
$ws.Range("F3").Value = 37563
$ws.Range("F6").Value = 767
$ws.Range("F7").Value = 475
$ws.Range("F11").Value = 384
$ws.Range("F15").Value = 841
$ws.Range("F17").Value = 708
$ws.Range("F19").Value = 34
$ws.Range("F25").Value = 14
$ws.Range("F26").Value = 9
$ws.Range("F27").Value = 646
$ws.Range("F31").Value = 1164
$ws.Range("F33").Value = 822
$ws.Range("F34").Value = 2514
$ws.Range("F35").Value = 997
$ws.Range("F36").Value = 561
$ws.Range("F37").Value = 106
$ws.Range("F38").Value = 1156
$ws.Range("F41").Value = 765
$ws.Range("F42").Value = 55
$ws.Range("F43").Value = 1152
